$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was logged. It belongs at row 185 (pushing the
# existing rows 185:286 down by one, to 186:287) rather than at the very end.
$ws.Rows(185).Insert()

$ws.Cells.Item(185, 1).Value = 10
$ws.Cells.Item(185, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(185, 3).Value = "La Araucanía"
$ws.Cells.Item(185, 4).Value = 44572
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 100112037
$ws.Cells.Item(185, 7).Value = "Cebollín"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 65
$ws.Cells.Item(185, 11).Value = 5000
$ws.Cells.Item(185, 12).Value = 5000
$ws.Cells.Item(185, 13).Value = 5000
$ws.Cells.Item(185, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(185, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(185, 16).Value = 417
$ws.Cells.Item(185, 17).Value = 12
$ws.Cells.Item(185, 18).Value = "Hortaliza"
